$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @(3, 1, 7.572976999999999, 22.718931, 0.4497670593913077, 0.4497670593913078, 3, 1, 1.307106666666667, 3.92132, 0.01256263154946851, 0.01256263154946851, 9.898688723213331, 89.08819850892, 0.005650257850220919, 0.00565025785022092),
    @(3, 1, 7.572976999999999, 22.718931, 0.4497670593913077, 0.4497670593913078, 3, 1, 80.22623699999998, 240.678711, 0.77105616682495, 0.77105616682495, 607.5514475975488, 5467.96302837794, 0.3467956647783914, 0.3467956647783914),
    @(3, 1, 7.572976999999999, 22.718931, 0.4497670593913077, 0.4497670593913078, 3, 1, 22.51385866666667, 67.541576, 0.2163812016255815, 0.2163812016255815, 170.4969338639173, 1534.472404775256, 0.09732113676269544, 0.09732113676269545),
    @(3, 1, 2.784025666666667, 8.352077, 0.1653462089435359, 0.1653462089435359, 3, 1, 1.307106666666667, 3.92132, 0.01256263154946851, 0.01256263154946851, 3.639018509071111, 32.75116658164, 0.002077183501059077, 0.002077183501059077),
    @(3, 1, 2.784025666666667, 8.352077, 0.1653462089435359, 0.1653462089435359, 3, 1, 80.22623699999998, 240.678711, 0.77105616682495, 0.77105616682495, 223.3519029480829, 2010.167126532747, 0.1274912140670401, 0.1274912140670401),
    @(3, 1, 2.784025666666667, 8.352077, 0.1653462089435359, 0.1653462089435359, 3, 1, 22.51385866666667, 67.541576, 0.2163812016255815, 0.2163812016255815, 62.67916038370578, 564.1124434533521, 0.03577781137543677, 0.03577781137543677),
    @(3, 1, 6.480551000000001, 19.441653, 0.3848867316651562, 0.3848867316651562, 3, 1, 1.307106666666667, 3.92132, 0.01256263154946851, 0.01256263154946851, 8.470771415773335, 76.23694274196001, 0.004835190198188512, 0.004835190198188512),
    @(3, 1, 6.480551000000001, 19.441653, 0.3848867316651562, 0.3848867316651562, 3, 1, 80.22623699999998, 240.678711, 0.77105616682495, 0.77105616682495, 519.910220416587, 4679.191983749283, 0.2967692879795185, 0.2967692879795185),
    @(3, 1, 6.480551000000001, 19.441653, 0.3848867316651562, 0.3848867316651562, 3, 1, 22.51385866666667, 67.541576, 0.2163812016255815, 0.2163812016255815, 145.9022092961254, 1313.119883665128, 0.08328225348744923, 0.08328225348744925)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $rowVals = $data[$i]
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $colNum = $j + 5
        $ws.Cells.Item($rowNum, $colNum).Value = $rowVals[$j]
    }
}
